# "Unlock Conversation Effect and Unlock Note Effect Update"
#
# Row 8 ("Effect3" / [101:3]) used to duplicate the "Unlock Note Entry"
# effect text from row 5/6. This changes it to a distinct new effect,
# "Update old note entry" (更新旧笔记条目), rendered with the same
# 微软雅黑 font/style already used by the "End Conversation" row (B7) -
# this also adds a new shared string and bumps the row height the same
# way row 7's taller font does.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New text for B8, using the same font as B7 ("结束对话") so it reuses
# the existing cell style instead of fabricating a new one.
$ws.Range("B8").Value = "更新旧笔记条目"
$ws.Range("B8").Font.Name = "微软雅黑"

# Row 8 grows to the same height as the other 微软雅黑 rows (4 and 7).
$ws.Rows.Item(8).RowHeight = 15.6

# Active selection moves to D7.
$ws.Range("D7").Select()
